$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.735.37"
$ws.Range("E2").Value = "  -4.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.217.46"
$ws.Range("E3").Value = "  -5.94%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.19"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.51"
$ws.Range("E6").Value = "  -7.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.591"
$ws.Range("E7").Value = "  -6.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.562"
$ws.Range("E9").Value = "  -7.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.11"
$ws.Range("E10").Value = "  -8.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.96"
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("E12").Value = "  -9.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.78"
$ws.Range("E13").Value = "  -7.11%  "
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.862"
$ws.Range("E15").Value = "  -11.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.555.56"
$ws.Range("E16").Value = "  -5.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.23"
$ws.Range("E17").Value = "  -6.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.220.98"
$ws.Range("E18").Value = "  -5.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.689.62"
$ws.Range("E19").Value = "  -4.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.03"
$ws.Range("E20").Value = "  +5.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -8.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.43"
$ws.Range("E22").Value = "  -10.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.35"
$ws.Range("E23").Value = "  -10.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.15"
$ws.Range("E24").Value = "  -9.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.59"
$ws.Range("E25").Value = "  -8.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("E26").Value = "  -6.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").Value = "  -8.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.38"
$ws.Range("E30").Value = "  -10.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0897"
$ws.Range("E31").Value = "  -6.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.43"
$ws.Range("E32").Value = "  -8.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.10"
$ws.Range("E33").Value = "  -7.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.41"
$ws.Range("E34").Value = "  -6.55%  "
$ws.Range("E35").Value = "  -6.28%  "
$ws.Range("E36").Value = "  +9.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.96"
$ws.Range("E37").Value = "  +13.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.122"
$ws.Range("E38").Value = "  -5.86%  "
$ws.Range("E39").Value = "  -5.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.89"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  -9.87%  "
$ws.Range("E42").Value = "  -7.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.927.96"
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.34"
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.58"
$ws.Range("E46").Value = "  -11.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.207"
$ws.Range("E47").Value = "  -8.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.37"
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "60.37"
$ws.Range("E49").Value = "  -12.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.29"
$ws.Range("E50").Value = "  -6.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.864"
$ws.Range("E51").Value = "  +17.52%  "
